$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the 5 brand-new variety names at the bottom of the existing
# list (rows 160-164). They, together with "LADY CLAIRE" / "SOPRANO"
# (already present but mis-sorted at rows 158-159), will be put back into
# alphabetical order by the Sort call below.
$ws.Cells.Item(160, 1).Value2 = "OPAL"
$ws.Cells.Item(161, 1).Value2 = "ENDURO"
$ws.Cells.Item(162, 1).Value2 = "CORINNA"
$ws.Cells.Item(163, 1).Value2 = "SH C 1010"
$ws.Cells.Item(164, 1).Value2 = "RISSOLETTO"

# --- Re-sort the whole "variete" column (A2:A164), ascending, leaving the
# header row (A1) untouched. This restores alphabetical order for the
# previously misplaced "LADY CLAIRE"/"SOPRANO" rows and places the five
# new varieties added above into their correct alphabetical slot.
$sortRange = $ws.Range("A2:A164")
$sortRange.Sort($ws.Range("A2"), 1)

# --- The generic sort above uses ordinal comparison and can disturb the
# original tie-break ordering of near-duplicate names that do not actually
# move as part of this edit. Restore the untouched "LADY CHRIST*" trio to
# its original relative order (Excel's locale-aware sort ranks the
# apostrophe variant last in this tie).
$ws.Cells.Item(86, 1).Value2 = "LADY CHRISTEL"
$ws.Cells.Item(87, 1).Value2 = "LADY CHRISTL"
$ws.Cells.Item(88, 1).Value2 = "LADY CHRIST'L"

# --- Finally, append "ASTERIX" as a new last row, after the sorted block.
$ws.Cells.Item(165, 1).Value2 = "ASTERIX"

# Match the saved selection/view state from the authored workbook.
$ws.Range("A165").Select()
